# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook
$missing = [System.Reflection.Missing]::Value

# --- Existing sheet (o_10) ---
$ws1 = $wb.Worksheets.Item(1)

# --- Add o_20 sheet after o_10 ---
$ws2 = $wb.Worksheets.Add($missing, $ws1)
$ws2.Name = "o_20"

# --- Add o_20_jumbled sheet after o_20 ---
$ws3 = $wb.Worksheets.Add($missing, $ws2)
$ws3.Name = "o_20_jumbled"

# --- Header texts (shared across sheets) ---
$hdrPrompt = @"
prompt
"@
$hdrSolution = @"
solution
"@
$hdrLlmResponse = @"
llm_response
"@
$hdrEvalResponse = @"
evaluator_response
"@
$hdrEvalPartial = @"
evaluator_partial_correctness
"@

# --- Copy header formatting (bold, border, centered) from A1 to E1 on each sheet ---

# --- Populate o_10 (sheet1) ---
$p = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node P? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P
 A 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0
 F 0 1 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 G 0 0 1 0 0 1 0 1 0 0 1 0 0 0 0 0
 H 0 0 0 1 0 0 1 0 0 0 0 1 0 0 0 0
 I 0 0 0 0 1 0 0 0 0 0 0 0 1 0 0 0
 J 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 K 0 0 0 0 0 0 1 0 0 0 0 0 0 0 1 0
 L 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0
 N 0 0 0 0 0 0 0 0 0 1 0 0 1 0 1 0
 O 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
    
"@
$sol = @"
A -> E -> I -> M -> N -> O -> P
"@
$llm = @"
The shortest path from node A to node P is A - E- I- M- N- O- P.
"@
$evr = @"
invalid input
"@
$evp = @"
7/7
"@
$ws1.Range("A1").Value2 = $hdrPrompt
$ws1.Range("B1").Value2 = $hdrSolution
$ws1.Range("C1").Value2 = $hdrLlmResponse
$ws1.Range("D1").Value2 = $hdrEvalResponse
$ws1.Range("E1").Value2 = $hdrEvalPartial
$ws1.Range("A2").Value2 = $p
$ws1.Range("B2").Value2 = $sol
$ws1.Range("C2").Value2 = $llm
$ws1.Range("D2").Value2 = $evr
$ws1.Range("E2").Value2 = $evp
$ws1.Rows.Item(2).AutoFit() | Out-Null

# --- Populate o_20 (sheet2) ---
$p = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0
    
"@
$sol = @"
A -> F -> G -> H -> M -> N -> O -> T -> Y
"@
$llm = @"
The shortest path from node A to node Y is: A - F - G - H - I - D - E - J - O - T - Y.
"@
$evr = @"
invalid input
"@
$evp = @"
3/9
"@
$ws2.Range("A1").Value2 = $hdrPrompt
$ws2.Range("B1").Value2 = $hdrSolution
$ws2.Range("C1").Value2 = $hdrLlmResponse
$ws2.Range("D1").Value2 = $hdrEvalResponse
$ws2.Range("E1").Value2 = $hdrEvalPartial
$ws2.Range("A2").Value2 = $p
$ws2.Range("B2").Value2 = $sol
$ws2.Range("C2").Value2 = $llm
$ws2.Range("D2").Value2 = $evr
$ws2.Range("E2").Value2 = $evp
$ws2.Rows.Item(2).AutoFit() | Out-Null

# --- Populate o_20_jumbled (sheet3) ---
$p = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0
    
"@
$sol = @"
A -> B -> C -> H -> M -> N -> S -> T -> Y
"@
$llm = @"
The shortest path from node A to node Y is: A -> F -> K -> P -> U -> V -> W -> X -> Y.
"@
$evr = @"
invalid input
"@
$evp = @"
1/9
"@
$ws3.Range("A1").Value2 = $hdrPrompt
$ws3.Range("B1").Value2 = $hdrSolution
$ws3.Range("C1").Value2 = $hdrLlmResponse
$ws3.Range("D1").Value2 = $hdrEvalResponse
$ws3.Range("E1").Value2 = $hdrEvalPartial
$ws3.Range("A2").Value2 = $p
$ws3.Range("B2").Value2 = $sol
$ws3.Range("C2").Value2 = $llm
$ws3.Range("D2").Value2 = $evr
$ws3.Range("E2").Value2 = $evp
$ws3.Rows.Item(2).AutoFit() | Out-Null

# --- Apply header style (copied from existing A1:D1 header format) to new E1 header cells ---
$ws1.Range("A1").Copy() | Out-Null
$ws1.Range("E1").PasteSpecial(-4122)
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A1:E1").PasteSpecial(-4122)
$ws1.Range("A1").Copy() | Out-Null
$ws3.Range("A1:E1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Keep o_10 as the active/selected sheet (matches tabSelected on sheet1 only) ---
$ws1.Activate()
